$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '36.435.82'
$ws.Range("E2").Value = '  +0.32%  '

# Row 3
$ws.Range("D3").Value = '2.010.06'
$ws.Range("E3").Value = '  -1.53%  '

# Row 4
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
Set-TextValue $ws.Range("D5") '252.13'
$ws.Range("E5").Value = '  +3.09%  '

# Row 6
$ws.Range("E6").Value = '  -2.96%  '

# Row 7
Set-TextValue $ws.Range("D7") '61.45'
$ws.Range("E7").Value = '  +12.82%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.370'
$ws.Range("E9").Value = '  +1.58%  '

# Row 10
Set-TextValue $ws.Range("D10") '58.49'
$ws.Range("E10").Value = '  -1.21%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0742'
$ws.Range("E11").Value = '  +0.30%  '

# Row 12
$ws.Range("E12").Value = '  -1.84%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.901'
$ws.Range("E13").Value = '  -0.49%  '

# Row 14
Set-TextValue $ws.Range("D14") '14.86'
$ws.Range("E14").Value = '  +3.88%  '

# Row 15
$ws.Range("D15").Value = '2.301.26'
$ws.Range("E15").Value = '  -1.66%  '

# Row 16
Set-TextValue $ws.Range("D16") '20.23'
$ws.Range("E16").Value = '  +15.83%  '

# Row 17
Set-TextValue $ws.Range("D17") '5.44'
$ws.Range("E17").Value = '  +1.90%  '

# Row 18
$ws.Range("D18").Value = '2.018.80'
$ws.Range("E18").Value = '  -1.08%  '

# Row 19
$ws.Range("D19").Value = '36.338.23'
$ws.Range("E19").Value = '  +0.34%  '

# Row 20
Set-TextValue $ws.Range("D20") '71.93'

# Row 21
$ws.Range("E21").Value = '  +0.89%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.25'
$ws.Range("E22").Value = '  +1.27%  '

# Row 23
Set-TextValue $ws.Range("D23") '234.27'
$ws.Range("E23").Value = '  -0.73%  '

# Row 24
$ws.Range("E24").Value = '  +21.87%  '

# Row 25
$ws.Range("E25").Value = '  -0.05%  '

# Row 26
$ws.Range("E26").Value = '  -1.17%  '

# Row 27
Set-TextValue $ws.Range("D27") '9.55'
$ws.Range("E27").Value = '  +2.54%  '

# Row 28
Set-TextValue $ws.Range("D28") '163.72'
$ws.Range("E28").Value = '  -0.16%  '

# Row 29
Set-TextValue $ws.Range("D29") '19.62'
$ws.Range("E29").Value = '  -1.30%  '

# Row 30
$ws.Range("E30").Value = '  -0.72%  '

# Row 31
Set-TextValue $ws.Range("D31") '5.11'
$ws.Range("E31").Value = '  +3.02%  '

# Row 32
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D32") '0.111'
$ws.Range("E32").Value = '  +23.35%  '

# Row 33
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D33") '1.18'
$ws.Range("E33").Value = '  +1.00%  '

# Row 34
Set-TextValue $ws.Range("D34") '4.56'
$ws.Range("E34").Value = '  +5.07%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0608'
$ws.Range("E35").Value = '  +1.65%  '

# Row 36
Set-TextValue $ws.Range("D36") '2.44'
$ws.Range("E36").Value = '  +10.52%  '

# Row 37
$ws.Range("E37").Value = '  -0.08%  '

# Row 38
$ws.Range("E38").Value = '  -0.82%  '

# Row 39
$ws.Range("E39").Value = '  +17.14%  '

# Row 40
$ws.Range("E40").Value = '  +14.14%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D41") '1.23'
$ws.Range("E41").Value = '  +1.60%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D42") '2.77'
$ws.Range("E42").Value = '  +22.39%  '

# Row 43
$ws.Range("E43").Value = '  +1.22%  '

# Row 44
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D44") '1.13'
$ws.Range("E44").Value = '  +3.03%  '

# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D45") '8.05'
$ws.Range("E45").Value = '  +8.03%  '

# Row 46
$ws.Range("E46").Value = '  +0.54%  '

# Row 47
Set-TextValue $ws.Range("D47") '16.79'
$ws.Range("E47").Value = '  +7.98%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D48") '94.43'
$ws.Range("E48").Value = '  +2.03%  '

# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.429.65'
$ws.Range("E49").Value = '  +2.09%  '

# Row 50
$ws.Range("E50").Value = '  -0.98%  '

# Row 51
Set-TextValue $ws.Range("D51") '46.97'
$ws.Range("E51").Value = '  +2.39%  '
